# Move the hidden "_GoBack" bookmark from the end of the "Email 7:" line
# (right after the table-of-emails block) down to the last (empty) cell of
# the final table in the document.
#
# In the source document, _GoBack sits collapsed right before the paragraph
# mark of the paragraph that ends with "Email 7: ". The edit relocates it to
# sit (also collapsed) right before the paragraph mark of the empty
# paragraph in the bottom-right cell of the last table ("My Analysis"
# table), which is where Word's cursor/_GoBack caret position ended up the
# next time the document was saved.

$d = $word.ActiveDocument
$bookmarks = $d.Bookmarks

if ($bookmarks.Exists("_GoBack")) {
    $old = $bookmarks.Item("_GoBack")
    $old.Delete()
}

# Locate the last table in the document and its bottom-right cell.
$lastTable = $d.Tables.Item($d.Tables.Count)
$lastRow = $lastTable.Rows.Count
$lastCol = $lastTable.Columns.Count
$targetCell = $lastTable.Cell($lastRow, $lastCol)
$targetParagraph = $targetCell.Range.Paragraphs.Item(1)

# Re-add the bookmark (using the paragraph's own Range, which — once added —
# Word collapses to sit right before the paragraph mark, matching the
# original collapsed placement pattern).
$bookmarks.Add("_GoBack", $targetParagraph.Range)
